$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update label text from "Cá nhân" to "Học tập"
$ws.Range("A2").Value = "Học tập"

# Update target hours value from 30 to 70
$ws.Range("C2").Value = 70

# Update the active selection to C2 (matches the saved view state)
$ws.Range("C2").Select()
